$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("id", "name", "status", "species", "type", "gender", "origin", "location", "image", "episode", "url", "created")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}
